$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.382.12'
$ws.Range("D3").Value = '3.170.95'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.54'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.57'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '3.168.64'
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("E9").Value = '  +1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.31'
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D15").Value = '3.694.14'
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '3.170.77'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '63.403.89'
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.56'
$ws.Range("E19").Value = '  -1.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.52'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("E22").Value = '  -2.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.64'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.12'
$ws.Range("E24").Value = '  -0.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.19'
$ws.Range("E25").Value = '  -2.03%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.69'
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.81'
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.73'
$ws.Range("E30").Value = '  -2.33%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.06'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.25'
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("E35").Value = '  -2.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.88'
$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("D37").Value = '0.0₃0737'
$ws.Range("E37").Value = '  +6.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '51.38'
$ws.Range("E38").Value = '  -0.36%  '
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.112'
$ws.Range("E41").Value = '  -2.61%  '
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '391.97'
$ws.Range("E43").Value = '  -5.20%  '
$ws.Range("D44").Value = '2.788.11'
$ws.Range("E44").Value = '  -7.22%  '
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '127.45'
$ws.Range("E46").Value = '  +2.20%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.75'
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("E49").Value = '  -2.93%  '
$ws.Range("E50").Value = '  -0.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.01'
$ws.Range("E51").Value = '  -3.42%  '
